$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Recomputed secant-method values for rows 2-10 (columns A-D).
$data = @(
    @(0, 1,                [double]"-2",                     [double]"1.000001"),
    @(0, 2,                [double]"4",                       [double]"1.000001"),
    @(1, 1.33333333333333, [double]"-0.962962962962963",      [double]"0.666666666666667"),
    @(2, 1.46268656716418, [double]"-0.333338874795104",      [double]"0.129353233830846"),
    @(3, 1.5311694321412,  [double]"0.0586264177094502",      [double]"0.06848286497702551"),
    @(4, 1.52092642051528, [double]"-0.0026933002019943",     [double]"0.0102430116259242"),
    @(5, 1.52137631666974, [double]"-2.01501923870355e-05",   [double]"0.0004498961544632"),
    @(6, 1.52137970798487, [double]"7.015478012664331e-09",   [double]"3.39131512783197e-06"),
    @(7, 1.52137970680456, [double]"-1.8429781545085e-14",    [double]"1.18030718532225e-09")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
}

# Remove the now-obsolete 11th row (the table shrank from 10 iterations to 9).
$ws.Rows("11").Delete()
